# fix extra time sum
# Remove the blank row between the header row and the "مجموع" (total) summary
# row on the first worksheet, so the summary row shifts up from row 4 to row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(2).Delete()
